$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.646.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.849.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  -0.60%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("E6").Value = "  -0.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4228"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.49%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3637"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.52"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07281"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8737"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.865.96"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.327"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.510"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06847"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "79.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008904"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.661.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.979"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("E24").Value = "  -4.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.091.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("E26").Value = "  -3.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "122.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +10.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.257"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.874"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +15.13%  "
$ws.Range("E32").Value = "  -0.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7649"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.73%  "
$ws.Range("E34").Value = "  +1.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.537"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.102"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.92%  "
$ws.Range("E37").Value = "  -0.77%  "
$ws.Range("E38").Value = "  +0.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05354"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01929"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.809"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.865"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5084"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.313"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06535"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.71%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4698"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.29"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.49%  "
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.621"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.35%  "
